$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 7.785765451129849
$ws.Cells.Item(2, 4).Value = 4.31035753193636
$ws.Cells.Item(2, 5).Value = 16.50157647476988
$ws.Cells.Item(2, 6).Value = 25.39530495257704
$ws.Cells.Item(2, 7).Value = 31.53344868302607
$ws.Cells.Item(2, 8).Value = 14.32394333831859
$ws.Cells.Item(2, 11).Value = 11.68655257041071
$ws.Cells.Item(2, 14).Value = 17.95006968125875
$ws.Cells.Item(3, 2).Value = 7.713313600657076
$ws.Cells.Item(3, 4).Value = 4.318801759986284
$ws.Cells.Item(3, 5).Value = 15.56285113366674
$ws.Cells.Item(3, 6).Value = 25.20059355215029
$ws.Cells.Item(3, 7).Value = 31.09416623933375
$ws.Cells.Item(3, 8).Value = 14.33171692962111
$ws.Cells.Item(3, 11).Value = 11.16617592998095
$ws.Cells.Item(3, 14).Value = 18.02170305331745
$ws.Cells.Item(4, 2).Value = 7.670353462034604
$ws.Cells.Item(4, 4).Value = 4.324181694981841
$ws.Cells.Item(4, 5).Value = 14.9619650266093
$ws.Cells.Item(4, 6).Value = 25.08960230092062
$ws.Cells.Item(4, 7).Value = 30.83519837125831
$ws.Cells.Item(4, 8).Value = 14.34021265029237
$ws.Cells.Item(4, 11).Value = 10.84420221602267
$ws.Cells.Item(4, 14).Value = 18.06759779777195
$ws.Cells.Item(5, 2).Value = 7.653249913651799
$ws.Cells.Item(5, 4).Value = 4.326423523078789
$ws.Cells.Item(5, 5).Value = 14.71120727542089
$ws.Cells.Item(5, 6).Value = 25.04656587895609
$ws.Cells.Item(5, 7).Value = 30.73250955627219
$ws.Cells.Item(5, 8).Value = 14.34460786340713
$ws.Cells.Item(5, 11).Value = 10.71077995491267
$ws.Cells.Item(5, 14).Value = 18.08678276729035
$ws.Cells.Item(6, 2).Value = 7.650434779130951
$ws.Cells.Item(6, 4).Value = 4.326798775812544
$ws.Cells.Item(6, 5).Value = 14.66922223883073
$ws.Cells.Item(6, 6).Value = 25.03955327435714
$ws.Cells.Item(6, 7).Value = 30.71563374547571
$ws.Cells.Item(6, 8).Value = 14.3453939482011
$ws.Cells.Item(6, 11).Value = 10.68849790072431
$ws.Cells.Item(6, 14).Value = 18.08999761105579
$ws.Cells.Item(7, 2).Value = 7.670121140519099
$ws.Cells.Item(7, 4).Value = 4.32421172831377
$ws.Cells.Item(7, 5).Value = 14.95860667923422
$ws.Cells.Item(7, 6).Value = 25.0890129661037
$ws.Cells.Item(7, 7).Value = 30.83380178685436
$ws.Cells.Item(7, 8).Value = 14.34026815181424
$ws.Cells.Item(7, 11).Value = 10.84241152074396
$ws.Cells.Item(7, 14).Value = 18.06785457705309
$ws.Cells.Item(8, 2).Value = 7.760479300878995
$ws.Cells.Item(8, 4).Value = 4.313228832259986
$ws.Cells.Item(8, 5).Value = 16.18313635031945
$ws.Cells.Item(8, 6).Value = 25.32641701658859
$ws.Cells.Item(8, 7).Value = 31.37983776590339
$ws.Cells.Item(8, 8).Value = 14.32584938279848
$ws.Cells.Item(8, 11).Value = 11.5016024265318
$ws.Cells.Item(8, 14).Value = 17.97437321880126
$ws.Cells.Item(9, 2).Value = 7.948888733360477
$ws.Cells.Item(9, 4).Value = 4.293222422821444
$ws.Cells.Item(9, 5).Value = 18.47025926389171
$ws.Cells.Item(9, 6).Value = 25.85790355745931
$ws.Cells.Item(9, 7).Value = 32.52949435554705
$ws.Cells.Item(9, 8).Value = 14.32722464697254
$ws.Cells.Item(9, 11).Value = 12.98184076344986
$ws.Cells.Item(9, 14).Value = 17.80613981030393
$ws.Cells.Item(10, 2).Value = 8.092924865315766
$ws.Cells.Item(10, 4).Value = 4.2794333634448
$ws.Cells.Item(10, 5).Value = 20.12035500672825
$ws.Cells.Item(10, 6).Value = 26.28560957210826
$ws.Cells.Item(10, 7).Value = 33.4132552035121
$ws.Cells.Item(10, 8).Value = 14.34643459583882
$ws.Cells.Item(10, 11).Value = 13.98438055717081
$ws.Cells.Item(10, 14).Value = 17.69161383268949
$ws.Cells.Item(11, 2).Value = 8.15940119849162
$ws.Cells.Item(11, 4).Value = 4.273352988582898
$ws.Cells.Item(11, 5).Value = 20.82910540226662
$ws.Cells.Item(11, 6).Value = 26.48755809502172
$ws.Cells.Item(11, 7).Value = 33.82183133405616
$ws.Cells.Item(11, 8).Value = 14.35914105179858
$ws.Cells.Item(11, 11).Value = 14.41580243235096
$ws.Cells.Item(11, 14).Value = 17.64145774673831
$ws.Cells.Item(12, 2).Value = 8.184688534835988
$ws.Cells.Item(12, 4).Value = 4.271077799848347
$ws.Cells.Item(12, 5).Value = 21.09150582798643
$ws.Cells.Item(12, 6).Value = 26.56502871381348
$ws.Cells.Item(12, 7).Value = 33.97732620113138
$ws.Cells.Item(12, 8).Value = 14.36452349344003
$ws.Cells.Item(12, 11).Value = 14.57561907280052
$ws.Cells.Item(12, 14).Value = 17.62274230672197
$ws.Cells.Item(13, 2).Value = 8.179237740143449
$ws.Cells.Item(13, 4).Value = 4.271566592813239
$ws.Cells.Item(13, 5).Value = 21.0352589586018
$ws.Cells.Item(13, 6).Value = 26.5483007343501
$ws.Cells.Item(13, 7).Value = 33.9438057058062
$ws.Cells.Item(13, 8).Value = 14.36333890514032
$ws.Cells.Item(13, 11).Value = 14.54135785594299
$ws.Cells.Item(13, 14).Value = 17.6267606935372
$ws.Cells.Item(14, 2).Value = 8.161479445401932
$ws.Cells.Item(14, 4).Value = 4.273165261876184
$ws.Cells.Item(14, 5).Value = 20.85081304975229
$ws.Cells.Item(14, 6).Value = 26.49391204358995
$ws.Cells.Item(14, 7).Value = 33.83460938135165
$ws.Cells.Item(14, 8).Value = 14.35957242994365
$ws.Cells.Item(14, 11).Value = 14.42902190278481
$ws.Cells.Item(14, 14).Value = 17.63991246296481
$ws.Cells.Item(15, 2).Value = 8.150616168754206
$ws.Cells.Item(15, 4).Value = 4.274148040406357
$ws.Cells.Item(15, 5).Value = 20.7370558376652
$ws.Cells.Item(15, 6).Value = 26.46072529830625
$ws.Cells.Item(15, 7).Value = 33.76781964568672
$ws.Cells.Item(15, 8).Value = 14.35733968470073
$ws.Cells.Item(15, 11).Value = 14.35974989568026
$ws.Cells.Item(15, 14).Value = 17.64800439956062
$ws.Cells.Item(16, 2).Value = 8.088597707730896
$ws.Cells.Item(16, 4).Value = 4.27983457121316
$ws.Cells.Item(16, 5).Value = 20.07319665725091
$ws.Cells.Item(16, 6).Value = 26.27255485703306
$ws.Cells.Item(16, 7).Value = 33.38667136179097
$ws.Cells.Item(16, 8).Value = 14.34568407269759
$ws.Cells.Item(16, 11).Value = 13.95568932276941
$ws.Cells.Item(16, 14).Value = 17.69493058827675
$ws.Cells.Item(17, 2).Value = 8.050779299303864
$ws.Cells.Item(17, 4).Value = 4.2833720903345
$ws.Cells.Item(17, 5).Value = 19.65523757679547
$ws.Cells.Item(17, 6).Value = 26.15896412815482
$ws.Cells.Item(17, 7).Value = 33.15440810419199
$ws.Cells.Item(17, 8).Value = 14.33955031091142
$ws.Cells.Item(17, 11).Value = 13.70149165046344
$ws.Cells.Item(17, 14).Value = 17.72421453859569
$ws.Cells.Item(18, 2).Value = 8.029118279905077
$ws.Cells.Item(18, 4).Value = 4.285424901222303
$ws.Cells.Item(18, 5).Value = 19.41089392284453
$ws.Cells.Item(18, 6).Value = 26.094328685766
$ws.Cells.Item(18, 7).Value = 33.02144269178243
$ws.Cells.Item(18, 8).Value = 14.33639582021857
$ws.Cells.Item(18, 11).Value = 13.55296504811015
$ws.Cells.Item(18, 14).Value = 17.74124081902894
$ws.Cells.Item(19, 2).Value = 8.021800581681003
$ws.Cells.Item(19, 4).Value = 4.286123070961987
$ws.Cells.Item(19, 5).Value = 19.32748478653457
$ws.Cells.Item(19, 6).Value = 26.07256620879228
$ws.Cells.Item(19, 7).Value = 32.97653552893514
$ws.Cells.Item(19, 8).Value = 14.33539189429292
$ws.Cells.Item(19, 11).Value = 13.50227863496509
$ws.Cells.Item(19, 14).Value = 17.7470370955378
$ws.Cells.Item(20, 2).Value = 8.054795868226927
$ws.Cells.Item(20, 4).Value = 4.282993642243079
$ws.Cells.Item(20, 5).Value = 19.70013803327323
$ws.Cells.Item(20, 6).Value = 26.17098417834836
$ws.Cells.Item(20, 7).Value = 33.17906931490344
$ws.Cells.Item(20, 8).Value = 14.34016460259319
$ws.Cells.Item(20, 11).Value = 13.72879150761469
$ws.Cells.Item(20, 14).Value = 17.72107829250193
$ws.Cells.Item(21, 2).Value = 8.166692570313968
$ws.Cells.Item(21, 4).Value = 4.272694955299484
$ws.Cells.Item(21, 5).Value = 20.90515151596225
$ws.Cells.Item(21, 6).Value = 26.50986078470504
$ws.Cells.Item(21, 7).Value = 33.86666324062036
$ws.Cells.Item(21, 8).Value = 14.36066324519166
$ws.Cells.Item(21, 11).Value = 14.46211419439014
$ws.Cells.Item(21, 14).Value = 17.63604194974683
$ws.Cells.Item(22, 2).Value = 8.240478569139542
$ws.Cells.Item(22, 4).Value = 4.266123242994943
$ws.Cells.Item(22, 5).Value = 21.6578112529046
$ws.Cells.Item(22, 6).Value = 26.73711672890758
$ws.Cells.Item(22, 7).Value = 34.32049330459522
$ws.Cells.Item(22, 8).Value = 14.37738694837132
$ws.Cells.Item(22, 11).Value = 14.92066909056191
$ws.Cells.Item(22, 14).Value = 17.58208294115437
$ws.Cells.Item(23, 2).Value = 8.201045208977011
$ws.Cells.Item(23, 4).Value = 4.269616244093679
$ws.Cells.Item(23, 5).Value = 21.25928308258693
$ws.Cells.Item(23, 6).Value = 26.61531889615125
$ws.Cells.Item(23, 7).Value = 34.07792286941216
$ws.Cells.Item(23, 8).Value = 14.36815687850211
$ws.Cells.Item(23, 11).Value = 14.67782726714495
$ws.Cells.Item(23, 14).Value = 17.61073447907764
$ws.Cells.Item(24, 2).Value = 8.052979722481385
$ws.Cells.Item(24, 4).Value = 4.283164679245586
$ws.Cells.Item(24, 5).Value = 19.67985116552336
$ws.Cells.Item(24, 6).Value = 26.16554782365762
$ws.Cells.Item(24, 7).Value = 33.1679182063371
$ws.Cells.Item(24, 8).Value = 14.3398857228941
$ws.Cells.Item(24, 11).Value = 13.71645667012516
$ws.Cells.Item(24, 14).Value = 17.72249559566321
$ws.Cells.Item(25, 2).Value = 7.896843147409571
$ws.Cells.Item(25, 4).Value = 4.298473267233454
$ws.Cells.Item(25, 5).Value = 17.82521834166479
$ws.Cells.Item(25, 6).Value = 25.70735522412978
$ws.Cells.Item(25, 7).Value = 32.21095098983222
$ws.Cells.Item(25, 8).Value = 14.32366341434336
$ws.Cells.Item(25, 11).Value = 12.59089941175386
$ws.Cells.Item(25, 14).Value = 17.85004906668011
